$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before Q (default_count, default_value) - this shifts
# the old Q..T (most_frequent_value, memory_consumed_bytes, pattern_count,
# patterns) to S..V.
$ws.Range("Q1:R1").EntireColumn.Insert()

# New header cells
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# Data rows 2..15: default_count = 0, default_value = "<Unspecified>" for every row
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 17).Value = 0              # Q = default_count
    $ws.Cells.Item($r, 18).Value = "<Unspecified>" # R = default_value
}

# The most_frequent_value column (now column S) changes for some rows.
$ws.Range("S3").Value = "Guilty"
$ws.Range("S4").Value = ""
$ws.Range("S6").Value = "Joseph"
$ws.Range("S7").Value = ""
$ws.Range("S8").Value = "Griffin"
$ws.Range("S9").Value = "19"
$ws.Range("S12").Value = "NORTHWESTERN"
$ws.Range("S13").Value = ""
$ws.Range("S14").Value = ""
$ws.Range("S15").Value = ""
